# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit
# to the NIT-9003913765 EC workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row before row 29 so the table grows from 13 to 14
#    period rows (rows 16-28 -> rows 16-29), and the footer rows (old 33/34)
#    shift down to 34/35 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("29").Insert()

# Copy the formatting (borders/fill/number format - the special "last row"
# look) that used to belong to row 28 down onto the freshly inserted row 29.
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Now give row 28 the "regular" (non-last-row) style that row 27 already
# uses, since row 28 is no longer the final row of the table.
$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Update the two summary fields that changed.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 773460      # VALOR MORA
$ws.Range("F13").Value = 14          # Cant. Periodos

# ---------------------------------------------------------------------------
# 3. Rewrite the worker/period detail table (rows 16-29).
# ---------------------------------------------------------------------------
$tipoDoc = "CC"

# Row 16 - Edwin Jovanys Diaz Zuñiga (unchanged)
$ws.Range("B16").Value = $tipoDoc
$ws.Range("C16").Value = "1047394420"
$ws.Range("D16").Value = "EDWIN JOVANYS DIAZ ZUÑIGA"
$ws.Range("E16").Value = "2109"
$ws.Range("F16").Value = 50776
$ws.Range("G16").Value = 2538780

# Row 17 - Jose Luis Cordoba Manrique (new, moved up)
$ws.Range("B17").Value = $tipoDoc
$ws.Range("C17").Value = "73184986"
$ws.Range("D17").Value = "JOSE LUIS CORDOBA MANRIQUE"
$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 7958
$ws.Range("G17").Value = 3202710

# Rows 18-20 - Rodolfo Sanjulian Trespalacios
$ws.Range("B18").Value = $tipoDoc
$ws.Range("C18").Value = "8850684"
$ws.Range("D18").Value = "RODOLFO SANJULIAN TRESPALACIOS"
$ws.Range("E18").Value = "2112"
$ws.Range("F18").Value = 72000
$ws.Range("G18").Value = 1800000

$ws.Range("B19").Value = $tipoDoc
$ws.Range("C19").Value = "8850684"
$ws.Range("D19").Value = "RODOLFO SANJULIAN TRESPALACIOS"
$ws.Range("E19").Value = "2201"
$ws.Range("F19").Value = 72000
$ws.Range("G19").Value = 1800000

$ws.Range("B20").Value = $tipoDoc
$ws.Range("C20").Value = "8850684"
$ws.Range("D20").Value = "RODOLFO SANJULIAN TRESPALACIOS"
$ws.Range("E20").Value = "2202"
$ws.Range("F20").Value = 48000
$ws.Range("G20").Value = 1800000

# Row 21 - Saul Navarro Diaz (unchanged)
$ws.Range("B21").Value = $tipoDoc
$ws.Range("C21").Value = "73107491"
$ws.Range("D21").Value = "SAUL NAVARRO DIAZ"
$ws.Range("E21").Value = "2203"
$ws.Range("F21").Value = 108962
$ws.Range("G21").Value = 2724060

# Rows 22-23 - Miguel Angel Hernandez Perez (new, moved up)
$ws.Range("B22").Value = $tipoDoc
$ws.Range("C22").Value = "1000222821"
$ws.Range("D22").Value = "MIGUEL ANGEL HERNANDEZ PEREZ"
$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 39858
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = $tipoDoc
$ws.Range("C23").Value = "1000222821"
$ws.Range("D23").Value = "MIGUEL ANGEL HERNANDEZ PEREZ"
$ws.Range("E23").Value = "2502"
$ws.Range("F23").Value = 32266
$ws.Range("G23").Value = 1423500

# Rows 24-29 - Diana Shirley Perez Amariles (periods 2503-2508, last one new)
$ws.Range("B24").Value = $tipoDoc
$ws.Range("C24").Value = "43653236"
$ws.Range("D24").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E24").Value = "2503"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

$ws.Range("B25").Value = $tipoDoc
$ws.Range("C25").Value = "43653236"
$ws.Range("D25").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E25").Value = "2504"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500

$ws.Range("B26").Value = $tipoDoc
$ws.Range("C26").Value = "43653236"
$ws.Range("D26").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E26").Value = "2505"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 1423500

$ws.Range("B27").Value = $tipoDoc
$ws.Range("C27").Value = "43653236"
$ws.Range("D27").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E27").Value = "2506"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500

$ws.Range("B28").Value = $tipoDoc
$ws.Range("C28").Value = "43653236"
$ws.Range("D28").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E28").Value = "2507"
$ws.Range("F28").Value = 56940
$ws.Range("G28").Value = 1423500

$ws.Range("B29").Value = $tipoDoc
$ws.Range("C29").Value = "43653236"
$ws.Range("D29").Value = "DIANA SHIRLEY PEREZ AMARILES"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 56940
$ws.Range("G29").Value = 1423500
